$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing formatting like trailing zeros).
# (Applied per-cell since multi-area comma ranges do not reliably propagate
# the NumberFormat change to every area in this runtime.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2: D2: 21.809.67 -> 21.830.46; E2: -1.38% -> -1.49%
$ws.Range("D2").Value = "21.830.46"
$ws.Range("E2").Value = "  -1.49%  "

# Row 3: D3: 1.540.61 -> 1.541.88; E3: -0.98% -> -1.08%
$ws.Range("D3").Value = "1.541.88"
$ws.Range("E3").Value = "  -1.08%  "

# Row 4: D4: 1.006 -> 1.005; E4: +0.64% -> +0.67%
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.67%  "

# Row 5: E5: +0.63% -> +0.62%
$ws.Range("E5").Value = "  +0.62%  "

# Row 6: D6: 289.50 -> 288.98; E6: +0.53% -> +0.26%
$ws.Range("D6").Value = "288.98"
$ws.Range("E6").Value = "  +0.26%  "

# Row 7: D7: 0.3942 -> 0.3923; E7: +3.97% -> +3.14%
$ws.Range("D7").Value = "0.3923"
$ws.Range("E7").Value = "  +3.14%  "

# Row 8: D8: 0.3198 -> 0.3195; E8: -2.78% -> -2.99%
$ws.Range("D8").Value = "0.3195"
$ws.Range("E8").Value = "  -2.99%  "

# Row 9: D9: 43.24 -> 43.16; E9: +0.09% -> -0.86%
$ws.Range("D9").Value = "43.16"
$ws.Range("E9").Value = "  -0.86%  "

# Row 10: D10: 0.07180 -> 0.07160; E10: -2.51% -> -3.05%
$ws.Range("D10").Value = "0.07160"
$ws.Range("E10").Value = "  -3.05%  "

# Row 11: D11: 1.063 -> 1.064; E11: -6.73% -> -7.13%
$ws.Range("D11").Value = "1.064"
$ws.Range("E11").Value = "  -7.13%  "

# Row 12: D12: 1.006 -> 1.005; E12: +0.63% -> +0.68%
$ws.Range("D12").Value = "1.005"
$ws.Range("E12").Value = "  +0.68%  "

# Row 13: D13: 5.628 -> 5.631; E13: -3.34% -> -3.42%
$ws.Range("D13").Value = "5.631"
$ws.Range("E13").Value = "  -3.42%  "

# Row 14: D14: 18.55 -> 18.57; E14: -8.13% -> -8.28%
$ws.Range("D14").Value = "18.57"
$ws.Range("E14").Value = "  -8.28%  "

# Row 15: D15: 6.626 -> 6.637; E15: -3.00% -> -3.09%
$ws.Range("D15").Value = "6.637"
$ws.Range("E15").Value = "  -3.09%  "

# Row 16: D16: 1.545.18 -> 1.549.46; E16: -1.21% -> +0.01%
$ws.Range("D16").Value = "1.549.46"
$ws.Range("E16").Value = "  +0.01%  "

# Row 17: E17: -0.40% -> -1.13%
$ws.Range("E17").Value = "  -1.13%  "

# Row 18: D18: 0.06554 -> 0.06560; E18: -0.94% -> -1.19%
$ws.Range("D18").Value = "0.06560"
$ws.Range("E18").Value = "  -1.19%  "

# Row 19: E19: -3.04% -> -3.27%
$ws.Range("E19").Value = "  -3.27%  "

# Row 20: E20: +0.54% -> +0.58%
$ws.Range("E20").Value = "  +0.58%  "

# Row 21: D21: 6.142 -> 6.144; E21: -3.95% -> -4.15%
$ws.Range("D21").Value = "6.144"
$ws.Range("E21").Value = "  -4.15%  "

# Row 22: D22: 15.28 -> 15.25; E22: -5.17% -> -5.47%
$ws.Range("D22").Value = "15.25"
$ws.Range("E22").Value = "  -5.47%  "

# Row 23: D23: 11.01 -> 11.02; E23: -6.02% -> -6.05%
$ws.Range("D23").Value = "11.02"
$ws.Range("E23").Value = "  -6.05%  "

# Row 24: D24: 2.382 -> 2.386; E24: +2.89% -> +3.09%
$ws.Range("D24").Value = "2.386"
$ws.Range("E24").Value = "  +3.09%  "

# Row 25: D25: 21.841.97 -> 21.842.79; E25: -1.26% -> -1.43%
$ws.Range("D25").Value = "21.842.79"
$ws.Range("E25").Value = "  -1.43%  "

# Row 26: D26: 2.374 -> 2.378; E26: -6.06% -> -6.32%
$ws.Range("D26").Value = "2.378"
$ws.Range("E26").Value = "  -6.32%  "

# Row 27: D27: 145.15 -> 145.16; E27: -3.61% -> -3.54%
$ws.Range("D27").Value = "145.16"
$ws.Range("E27").Value = "  -3.54%  "

# Row 28: D28: 18.37 -> 18.47; E28: -3.98% -> -3.79%
$ws.Range("D28").Value = "18.47"
$ws.Range("E28").Value = "  -3.79%  "

# Row 29: D29: 4.857 -> 4.860; E29: -1.12% -> -1.17%
$ws.Range("D29").Value = "4.860"
$ws.Range("E29").Value = "  -1.17%  "

# Row 30: D30: 1.717.01 -> 1.722.47; E30: -0.88% -> -0.24%
$ws.Range("D30").Value = "1.722.47"
$ws.Range("E30").Value = "  -0.24%  "

# Row 31: D31: 117.20 -> 117.15; E31: -3.72% -> -3.84%
$ws.Range("D31").Value = "117.15"
$ws.Range("E31").Value = "  -3.84%  "

# Row 32: D32: 0.9678 -> 0.9664; E32: -10.55% -> -10.68%
$ws.Range("D32").Value = "0.9664"
$ws.Range("E32").Value = "  -10.68%  "

# Row 33: D33: 5.877 -> 5.869; E33: -1.52% -> -1.64%
$ws.Range("D33").Value = "5.869"
$ws.Range("E33").Value = "  -1.64%  "

# Row 34: D34: 0.08232 -> 0.08229; E34: +0.07% -> -0.10%
$ws.Range("D34").Value = "0.08229"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35: D35: 8.965 -> 8.992; E35: -4.01% -> -3.89%
$ws.Range("D35").Value = "8.992"
$ws.Range("E35").Value = "  -3.89%  "

# Row 36: D36: 1.556 -> 1.543; E36: -15.93% -> -16.64%
$ws.Range("D36").Value = "1.543"
$ws.Range("E36").Value = "  -16.64%  "

# Row 37: D37: 0.06060 -> 0.06083; E37: -3.19% -> -2.95%
$ws.Range("D37").Value = "0.06083"
$ws.Range("E37").Value = "  -2.95%  "

# Row 38: D38: 0.02235 -> 0.02234; E38: -4.34% -> -4.56%
$ws.Range("D38").Value = "0.02234"
$ws.Range("E38").Value = "  -4.56%  "

# Row 39: D39: 5.109 -> 5.102; E39: -3.58% -> -3.76%
$ws.Range("D39").Value = "5.102"
$ws.Range("E39").Value = "  -3.76%  "

# Row 40: D40: 1.196 -> 1.198; E40: -4.68% -> -4.55%
$ws.Range("D40").Value = "1.198"
$ws.Range("E40").Value = "  -4.55%  "

# Row 41: D41: 0.2030 -> 0.2034; E41: -6.13% -> -6.21%
$ws.Range("D41").Value = "0.2034"
$ws.Range("E41").Value = "  -6.21%  "

# Row 43: E43: -3.91% -> -4.10%
$ws.Range("E43").Value = "  -4.10%  "

# Row 44: D44: 0.5761 -> 0.5756; E44: -4.97% -> -5.53%
$ws.Range("D44").Value = "0.5756"
$ws.Range("E44").Value = "  -5.53%  "

# Row 45: D45: 3.757 -> 3.748; E45: +0.50% -> +0.18%
$ws.Range("D45").Value = "3.748"
$ws.Range("E45").Value = "  +0.18%  "

# Row 46: D46: 12.94 -> 12.92; E46: -5.50% -> -5.97%
$ws.Range("D46").Value = "12.92"
$ws.Range("E46").Value = "  -5.97%  "

# Row 47: D47: 0.5542 -> 0.5543; E47: -5.36% -> -5.72%
$ws.Range("D47").Value = "0.5543"
$ws.Range("E47").Value = "  -5.72%  "

# Row 48: D48: 116.95 -> 116.92; E48: -4.48% -> -4.45%
$ws.Range("D48").Value = "116.92"
$ws.Range("E48").Value = "  -4.45%  "

# Row 49: D49: 1.860 -> 1.859; E49: -6.63% -> -6.95%
$ws.Range("D49").Value = "1.859"
$ws.Range("E49").Value = "  -6.95%  "

# Row 50: D50: 1.130 -> 1.128; E50: -4.01% -> -4.25%
$ws.Range("D50").Value = "1.128"
$ws.Range("E50").Value = "  -4.25%  "

# Row 51: D51: 0.06747 -> 0.06751; E51: -3.91% -> -3.96%
$ws.Range("D51").Value = "0.06751"
$ws.Range("E51").Value = "  -3.96%  "
